$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price text values are kept as text, not converted to numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated values from the crypto data refresh
$ws.Range("D2").Value = '66.831.14'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '3.068.81'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '575.38'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '168.39'
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '3.066.12'
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("D10").Value = '6.39'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  -3.49%  '
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("E14").Value = '  -4.28%  '
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("D16").Value = '3.580.98'
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").Value = '66.755.35'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").Value = '16.78'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").Value = '3.064.51'
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("D21").Value = '491.34'
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("E23").Value = '  -3.38%  '
$ws.Range("D24").Value = '82.81'
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").Value = '12.68'
$ws.Range("E25").Value = '  -6.10%  '
$ws.Range("E26").Value = '  -4.59%  '
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("D29").Value = '7.76'
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("E30").Value = '  -5.21%  '
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("D32").Value = '27.47'
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("D33").Value = '0.112'
$ws.Range("E33").Value = '  -3.51%  '
$ws.Range("E34").Value = '  -3.32%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").Value = '0.946'
$ws.Range("E36").Value = '  -3.19%  '
$ws.Range("E37").Value = '  -4.57%  '
$ws.Range("D38").Value = '46.54'
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E41").Value = '  -3.58%  '
$ws.Range("D42").Value = '8.30'
$ws.Range("E42").Value = '  -4.55%  '
$ws.Range("D43").Value = '2.751.49'
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("E44").Value = '  -3.14%  '
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").Value = '367.22'
$ws.Range("E46").Value = '  -3.44%  '
$ws.Range("D49").Value = '24.41'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("E51").Value = '  -1.88%  '
